# "added 4wk low sales check"
# Update the Forecast Comparison sheet (MyForecast, Inventory Coverage,
# Reorder Urgency, Seasonality Index) and the Summary sheet (forecast
# totals / min forecast) to reflect the new 4-week low-sales check logic.

$wb = $excel.ActiveWorkbook

# ---- Sheet: Forecast Comparison ----
$fc = $wb.Worksheets.Item("Forecast Comparison")

# Row 2 (W10)
$fc.Range("H2").Value = 2.96
$fc.Range("L2").Value = 0.85

# Row 3 (W11)
$fc.Range("D3").Value = 13
$fc.Range("H3").Value = 1.96
$fc.Range("L3").Value = 1.18

# Row 4 (W12)
$fc.Range("D4").Value = 13
$fc.Range("H4").Value = 0.98
$fc.Range("J4").Value = "Urgent"
$fc.Range("L4").Value = 0.89

# Row 5 (W13)
$fc.Range("D5").Value = 12
$fc.Range("H5").Value = 0
$fc.Range("L5").Value = 0.93

# Row 6 (W14)
$fc.Range("D6").Value = 12
$fc.Range("L6").Value = 0.96

# Row 7 (W15)
$fc.Range("L7").Value = 1.17

# Row 8 (W16)
$fc.Range("D8").Value = 12
$fc.Range("L8").Value = 1.06

# Row 9 (W17)
$fc.Range("D9").Value = 11
$fc.Range("L9").Value = 0.8100000000000001

# Row 10 (W18)
$fc.Range("D10").Value = 11
$fc.Range("L10").Value = 0.97

# Row 11 (W19)
$fc.Range("D11").Value = 11
$fc.Range("L11").Value = 1.03

# Row 12 (W20)
$fc.Range("D12").Value = 11
$fc.Range("L12").Value = 1.03

# Row 13 (W21)
$fc.Range("D13").Value = 10
$fc.Range("L13").Value = 1.14

# Row 14 (W22)
$fc.Range("D14").Value = 10
$fc.Range("L14").Value = 1.02

# Row 15 (W23)
$fc.Range("D15").Value = 10
$fc.Range("L15").Value = 0.84

# Row 16 (W24)
$fc.Range("D16").Value = 10
$fc.Range("L16").Value = 0.8100000000000001

# Row 17 (W25)
$fc.Range("D17").Value = 10
$fc.Range("L17").Value = 0.98

# ---- Sheet: Summary ----
# These "Value" cells are stored as text (not numbers) in the workbook,
# so prefix with an apostrophe to keep them text-typed after the edit
# (otherwise Excel auto-converts the numeric-looking string to a number).
$sm = $wb.Worksheets.Item("Summary")

$sm.Range("B9").Value = "'189"
$sm.Range("B10").Value = "'102"
$sm.Range("B11").Value = "'53"
$sm.Range("B12").Value = "'14"
$sm.Range("B14").Value = "'10"
